# Edit: slide 21 ("Examples"), the Content Placeholder shape, 3rd paragraph.
# The single run
#   online_resource = 856 ? (ind2 = "0" || ind2 = "1"), map(".*=>Yes
# gets the text changed to
#   online_resource = 856u ? (ind2 = "0" || ind2 = "1"), map(".*=>Yes
# and split into four separate runs:
#   "online_resource "
#   "= "
#   "856u "
#   "? (ind2 = \"0\" || ind2 = \"1\"), map(\".*=>Yes"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$shape = $s.Shapes.Item(2)
$para = $shape.TextFrame.TextRange.Paragraphs(3)

# Step 1: turn "856" into "856u" (chars 19-21 are "856").
$numPart = $para.Characters(19, 3)
$numPart.Text = "856u"

# Step 2: split off "online_resource " as its own run (chars 1-16).
$run1 = $para.Characters(1, 16)
$run1.Text = "online_resource "

# Step 3: split off "= " as its own run (chars 17-18).
$run2 = $para.Characters(17, 2)
$run2.Text = "= "

# Step 4: split off "856u " as its own run (chars 19-23, now that the
# "u" has been inserted).
$run3 = $para.Characters(19, 5)
$run3.Text = "856u "

# What's left (from char 24 to the end) remains as the final run with the
# unchanged text "? (ind2 = "0" || ind2 = "1"), map(".*=>Yes".

Write-Host "Final paragraph text:" $para.Text
